$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 192 - shifts existing rows 192:280 down to 193:281,
# carrying their values/styles along automatically.
$ws.Rows.Item(192).Insert()

# Populate the newly-inserted row 192 with the new record.
$ws.Range("A192").Value = 7
$ws.Range("B192").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C192").Value = "Ñuble"
$ws.Range("D192").Value = 44917
$ws.Range("E192").Value = 16
$ws.Range("F192").Value = 100112043
$ws.Range("G192").Value = "Pepino ensalada"
$ws.Range("H192").Value = "Sin especificar"
$ws.Range("I192").Value = "Primera"
$ws.Range("J192").Value = 120
$ws.Range("K192").Value = 13000
$ws.Range("L192").Value = 14000
$ws.Range("M192").Value = 13500
$ws.Range("N192").Value = "$/caja 80 unidades"
$ws.Range("O192").Value = "Región del Maule"
$ws.Range("P192").Value = 169
$ws.Range("Q192").Value = 80
$ws.Range("R192").Value = "Hortaliza"

# D192 needs the same datetime style as the rest of column D (numFmtId 165).
$ws.Range("D192").NumberFormat = $ws.Range("D193").NumberFormat
